$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell value updates from the upstream "cryptos list" data refresh.
# Column D holds numeric-looking price strings that must stay plain text
# (to preserve formats like "71.171.39", "1.00", "0.000323"), so those
# cells are forced to text format before/after the write. Columns B, C, E
# are non-numeric-looking strings and do not need that treatment.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '71.171.39'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.88%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.806.67'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.90%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '699.24'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +10.77%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '172.97'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.64%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.806.12'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.94%  '
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('E9').Value = '  +0.88%  '
$ws.Range('E10').Value = '  +2.82%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.53'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +11.41%  '
$ws.Range('E12').Value = '  +0.50%  '
$ws.Range('E13').Value = '  +8.15%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.26'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.33%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.452.97'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.98%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.806.76'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.98%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '71.165.67'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.85%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '17.87'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.22'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.90%  '
$ws.Range('E20').Value = '  +1.17%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.20'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +17.62%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '481.37'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.98%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.715'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.35%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '84.00'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.77%  '
$ws.Range('E25').Value = '  +0.23%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.36'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.35%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.17'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.37%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.50'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.16%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.958.37'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E30').Value = '  -0.02%  '
$ws.Range('E31').Value = '  +14.57%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.31'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.02%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.55'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +6.70%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '29.55'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.73%  '
$ws.Range('E35').Value = '  -1.32%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '9.22'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.11%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.11%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.757.44'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.84%  '
$ws.Range('E39').Value = '  +1.50%  '
$ws.Range('E40').Value = '  +5.57%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.98'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.11%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.21'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +11.10%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.970'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.95%  '
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.00'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('B45').Value = 'FLOKI'
$ws.Range('C45').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.000323'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +21.46%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '45.49'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.76%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '49.34'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.67%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '160.40'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.42%  '
$ws.Range('E50').Value = '  -1.44%  '
